# Adds a new row of data (Balicek 5 kg) to the offer list sheet,
# matching the layout/appearance of the other rows, then adjusts
# the active selection and autofits column D to the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 4 into row 5 (keeps the same cell styles as the rows above)
$ws.Rows("4").Copy()
$ws.Rows("5").Insert(-4121)  # xlShiftDown
$excel.CutCopyMode = 0

# New row 5 content (order matches how the values were originally entered)
$ws.Range("B5").Value = "2501 kc"
$ws.Range("C5").Value = "21.1.2023"
$ws.Range("A5").Value = "Balicek 5 kg"
$ws.Range("D5").Value = "Jazyk + Jitrnice"

# Column D needs to widen slightly to fit the new, longer description text
$ws.Columns.Item(4).ColumnWidth = 11.75

# Update the active selection to reflect the new last cell
$ws.Range("D5").Select()
